$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump the Date value (B8, next to the "Date" label in A8) ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value2 = "2024-03-19T13:17:15+00:00"

# --- Elements sheet: swap the two "Mapping" columns (AK <-> AL) ---
$el = $wb.Worksheets.Item("Elements")

# Swap header text (AK1 / AL1): "Mapping: RIM Mapping" <-> "Mapping: Spécification..."
$akHeader = $el.Range("AK1").Value2
$alHeader = $el.Range("AL1").Value2
$el.Range("AK1").Value2 = $alHeader
$el.Range("AL1").Value2 = $akHeader

# Swap only the data rows that actually carry a mapping value in AK or AL
# (rows 3, 5 and 6 - rows 2 and 4 are blank in both columns, so leave
# them exactly as-is rather than rewriting them to a different blank
# representation).
$rowsToSwap = @(3, 5, 6)
foreach ($r in $rowsToSwap) {
    $akCell = $el.Cells.Item($r, 37)
    $alCell = $el.Cells.Item($r, 38)
    $akVal = $akCell.Value2
    $alVal = $alCell.Value2
    $akCell.Value2 = $alVal
    $alCell.Value2 = $akVal
}

# Swap the column widths of AK (col 37) and AL (col 38) to match the
# content that now lives in each column.
$el.Columns.Item(37).ColumnWidth = 85.33333333333334
$el.Columns.Item(38).ColumnWidth = 24.166666666666664
